$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$variaveis = @(
    "variavel",
    "rendimento médio real habitual de todos os trabalhos(r`$/mês)",
    "rendimento médio real habitual do trabalho principal(r`$/mês)",
    "rendimento-hora médio real habitual do trabalho principal(r`$/hora)",
    "rendimento-hora médio real habitual de todos os trabalhos(r`$/hora)",
    "taxa de formalização",
    "índice gini",
    "domicílio cedido por familiar",
    "outra forma",
    "15 a 29 anos",
    "população em idade de trabalhar",
    "60 anos ou mais",
    "população desocupada",
    "proporção de pessoas com acesso simultâneo aos três serviços de saneamento básico(%)",
    "população ocupada em trabalhos formais",
    "total pessoas ocupadas(1 000 pessoas)",
    "população na força de trabalho",
    "ensino médio completo ou superior incompleto",
    "45 a 59 anos",
    "população ocupada",
    "taxa composta de subutilização",
    "população",
    "número de beneficiários de plano de saúde",
    "sem instrução ou fundamental incompleto",
    "taxa de desocupação",
    "domicílio alugado",
    "população na força de trabalho potencial",
    "taxa total mortalidade",
    "população subutilizada",
    "taxa de participação",
    "70 anos ou mais",
    "0 a 14 anos",
    "taxa de analfabetismo",
    "nível de ocupação",
    "ensino fundamental completo ou médio incompleto",
    "ensino superior completo",
    "saneamento basico total(1 000 pessoas)",
    "número mensal médio de leitos de internação (total)",
    "60 a 69 anos",
    "30 a 44 anos",
    "total pessoas por condição de ocupação a domicílio(1 000 pessoas)",
    "domicílio próprio - pagando",
    "domicílio próprio - já pago",
    "domicílio cedido por empregador",
    "domicílio cedido de outra forma"
)

$impactos = @(
    "impacto",
    0.5248516201972961,
    0.2970086336135864,
    0.1115423515439034,
    0.03483881801366806,
    0.01165385730564594,
    0.01049863174557686,
    0.003960499539971352,
    0.002395674120634794,
    0.001086007221601903,
    0.0006411526119336486,
    0.0002653948613442481,
    0.0001686291070654988,
    0.0001430262927897274,
    0.0001381486508762464,
    0.0001183039057650603,
    [double]"9.130220860242844e-05",
    [double]"8.127070759655908e-05",
    [double]"7.304239261429757e-05",
    [double]"6.763084093108773e-05",
    [double]"6.698773358948529e-05",
    [double]"4.928356793243438e-05",
    [double]"4.178608287475072e-05",
    [double]"4.1456580220256e-05",
    [double]"3.503988773445599e-05",
    [double]"3.201008803443983e-05",
    [double]"2.731939639488701e-05",
    [double]"2.432652217976283e-05",
    [double]"2.329562448721845e-05",
    [double]"2.285359005327336e-05",
    [double]"1.167465870821616e-05",
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0
)

for ($i = 0; $i -lt $variaveis.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $variaveis[$i]
    $ws.Cells.Item($row, 2).Value = $impactos[$i]
}

# B1 mirrors A1's header style (bold/centered/bordered) without minting a new style.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
